# إضافة حدث جديد في Card7 by admin at 2026-02-18 13:11:42
#
# This change adds a new maintenance/lubrication event row to the "Card7"
# worksheet (the event that already existed as the last row of "Card6"
# - "تشحيم" / lubrication, dated 18/2/2026, by "يوسف وابراهيم").
#
# Because the workbook is produced by a pandas/openpyxl export pipeline,
# re-exporting the data also toggles how "missing" values are rendered in
# the two affected sheets:
#   - Card7 : previously-blank cells in the data rows become the literal
#             text "nan" (the sheet's dimension grows from A1:O14 to
#             A1:O15 to include the newly appended row).
#   - Card6 : the event row that was moved out is removed (dimension
#             shrinks from A1:O14 to A1:O13) and its "nan" placeholder
#             text reverts back to blank cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Card7 : fill previously blank cells with the literal text "nan"
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Card7")

$card7NanCells = @(
  "D2","E2","F2","G2","H2","I2","J2","K2","N2",
  "G3","H3","I3","J3","K3",
  "D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4","O4",
  "D5","E5","F5","G5","H5","I5","J5","K5","L5","M5","N5","O5",
  "D6","E6","F6","G6","H6","I6","J6","K6","L6","M6","N6","O6",
  "D7","E7","F7","G7","H7","I7","J7","K7","L7","M7","N7","O7",
  "D8","E8","F8","G8","H8","I8","J8","K8","L8","M8","N8","O8",
  "D9","E9","F9","G9","H9","I9","J9","K9","L9","M9","N9","O9",
  "D10","E10","F10","G10","H10","I10","J10","K10","L10","M10","N10","O10",
  "D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11","O11",
  "D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12","O12",
  "B13","C13","D13","E13","F13","G13","H13","I13","J13","K13",
  "B14","C14","D14","E14","F14","G14","H14","I14","J14","K14"
)

foreach ($ref in $card7NanCells) {
  $ws7.Range($ref).Value = "nan"
}

# ---------------------------------------------------------------------
# 2) Card7 : append the new event as row 15
# ---------------------------------------------------------------------
$a15 = $ws7.Range("A15")
$a15.NumberFormat = "@"
$a15.Value = "7"

$ws7.Range("L15").Value = "18/2/2026"
$ws7.Range("M15").Value = "تشحيم"
$ws7.Range("N15").Value = "تم تشحيم السلندر  15جرام كل جانب   والدوفر 7 جرام كل جانب"
$ws7.Range("O15").Value = "يوسف وابراهيم"

# ---------------------------------------------------------------------
# 3) Card6 : clear the "nan" placeholders back to blank cells
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Card6")

$card6NanCells = @(
  "D2","E2","F2","G2","H2","I2","J2","K2","L2","M2","N2","O2",
  "G3","H3","I3","J3","K3",
  "D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4","O4",
  "D5","E5","F5","G5","H5","I5","J5","K5","L5","M5","N5","O5",
  "D6","E6","F6","G6","H6","I6","J6","K6","L6","M6","N6","O6",
  "D7","E7","F7","G7","H7","I7","J7","K7","L7","M7","N7","O7",
  "D8","E8","F8","G8","H8","I8","J8","K8","L8","M8","N8","O8",
  "D9","E9","F9","G9","H9","I9","J9","K9","L9","M9","N9","O9",
  "D10","E10","F10","G10","H10","I10","J10","K10","L10","M10","N10","O10",
  "D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11","O11",
  "D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12","O12",
  "B13","C13","D13","E13","F13","G13","H13","I13","J13","K13"
)

foreach ($ref in $card6NanCells) {
  $ws6.Range($ref).Value = ""
}

# ---------------------------------------------------------------------
# 4) Card6 : remove row 14 (the event that moved to Card7)
# ---------------------------------------------------------------------
$ws6.Rows.Item(14).Delete()
